$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns S..Y hold HYPERLINK(...) formulas that currently only take one
# argument (the URL). This update adds a second argument - the label text -
# which is the value found in column A of the same row (e.g. "A 10563-2022").
$linkCols = 19,20,21,22,23,24,25   # S,T,U,V,W,X,Y

$lastRow = 100
for ($r = 2; $r -le $lastRow; $r++) {

    # --- 1) Update the "Forandrad" (changed) date in column C ---
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Value2 -ne $null) {
        $cCell.Value2 = 45186
    }

    # --- 2) Add the display-text argument to every HYPERLINK formula ---
    $label = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -like 'HYPERLINK(*' -or $f -like '=HYPERLINK(*') {
                if ($f -notmatch ',\s*"') {
                    $idx = $f.LastIndexOf(")")
                    if ($idx -ge 0) {
                        $newF = $f.Substring(0, $idx) + ', "' + $label + '"' + $f.Substring($idx)
                        $cell.Formula = $newF
                    }
                }
            }
        }
    }
}
